$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Final target data for the parts list (rows 3-36, columns A/B/C).
# Row 2 stays blank (matches the original layout gap between header & data).
# ---------------------------------------------------------------------------
$rows = @(
    @{Row=3;  A='1276-1096-1-ND';          B='C2, C6';       C=2},
    @{Row=4;  A='1276-1286-1-ND';          B='C3,C4';        C=2},
    @{Row=5;  A='445-13454-1-ND';          B='C5';           C=1},
    @{Row=6;  A='PCF1412CT-ND';            B='C7';           C=1},
    @{Row=7;  A='ES2F-E3/52TGICT-ND';      B='D1';           C=1},
    @{Row=8;  A='IRF644SPBF-ND';           B='Q1';           C=1},
    @{Row=9;  A='MMBTA92LT1GOSCT-ND';      B='Q2-4';         C=3},
    @{Row=10; A='P50MCT-ND';               B='RSENSE';       C=1},
    @{Row=11; A='PT1.5MXCT-ND';            B='R1';           C=1},
    @{Row=12; A='490-6996-1-ND';           B='R2';           C=1},
    @{Row=13; A='541-3976-1-ND';           B='R3';           C=1},
    @{Row=14; A='RMCF0805FT10K0CT-ND';     B='R4,R6';        C=2},
    @{Row=15; A='RMCF0805FT300KCT-ND';     B='R5,R7';        C=2},
    @{Row=16; A='RNCF0805DTE10K0CT-ND';    B='R12,R13,R14';  C=3},
    @{Row=17; A='A130138CT-ND';            B='R10,R11';      C=2},
    @{Row=18; A='732-1248-1-ND';           B='L1';           C=1},
    @{Row=19; A='952-1745-1-ND';           B='BT1';          C=1},
    @{Row=20; A='MAX1771CSA+CT-ND';        B='IC1';          C=1},
    @{Row=21; A='296-12010-1-ND';          B='IC2';          C=1},
    @{Row=22; A='HV5122PG-G-ND';           B='IC3';          C=1},
    @{Row=23; A='AZ1117CH-3.3TRG1DICT-ND'; B='VR1';          C=1},
    @{Row=24; A='945-1648-5-ND';           B='VR2';          C=1},
    @{Row=25; A='CP-037A-ND';              B='IC5';          C=1},
    @{Row=26; A='DS3231MZ+-ND';            B='IC6';          C=1},
    @{Row=27; A='1528-1438-ND';            B='IC7';          C=1},
    @{Row=28; A='PTS645SM43SMTR92 LFS';    B='B1';           C=3},
    @{Row=29; A='S7121-ND';                B='P1';           C=1},
    @{Row=30; A='732-5317-ND';             B='P2';           C=1},
    @{Row=31; A='IRF644SPBF-ND';           B='Q1';           C=1},
    @{Row=32; A='MMBTA92LT1GOSCT-ND';      B='Q2,Q3';        C=2},
    @{Row=33; A='1276-1286-1-ND';          B='C3, C4';       C=2},
    @{Row=34; A='445-13454-1-ND';          B='C5';           C=1},
    @{Row=35; A='1276-1096-1-ND';          B='C2,C6';        C=2},
    @{Row=36; A='PCF1412CT-ND';            B='C7';           C=1}
)

# Rows whose column-A cell must carry a particular direct-format style
# (index into the existing cellXfs table). Stable, never-recoloured donor
# cells are used as the copy source so no new fonts/fills get minted.
#   style 1 -> existing A5 (Arial 9 #444444)
#   style 2 -> existing A14 (Arial 9 #444444, white fill, wrap)
#   style 3 -> existing A8 (Arial 9 black)
$styleDonor = @{1='A5'; 2='A14'; 3='A8'}

$styledRows = @{
    3  = 3
    5  = 1
    8  = 3
    9  = 1
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    19 = 2
    21 = 1
    22 = 1
    23 = 1
    26 = 1
    27 = 2
}

# Rows whose A cell currently owns a direct-format style that must be
# stripped back to the default (no explicit style) in the final layout.
$unstyledRows = @(18, 20, 24, 25)

# Apply/clear the column-A styles first (independent of the values, which
# get overwritten below regardless of what the donor copy brought along).
foreach ($r in $unstyledRows) {
    $ws.Range("A$r").ClearFormats()
}
foreach ($r in $styledRows.Keys) {
    $styleId = $styledRows[$r]
    $donor = $styleDonor[$styleId]
    $ws.Range($donor).Copy($ws.Range("A$r"))
}

# The little hyperlink-style placeholder cell in column F slides from F20
# down to F21/F22/F23. Clone its format onto the new spots, then drop it
# from its old one.
$ws.Range("F20").Copy($ws.Range("F21"))
$ws.Range("F20").Copy($ws.Range("F22"))
$ws.Range("F20").Copy($ws.Range("F23"))
$ws.Range("F20").Clear()

# Write every row's actual values (this also overwrites whatever value a
# style-donor Copy() above incidentally carried into column A).
foreach ($row in $rows) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
}

# Selection / scroll position, matching the saved view in the edited file.
$ws.Range("A10").Select()
$ws.Range("H21").Select()
